# Auto-generated Excel COM-interop script
# Applies scheduled-runner value updates across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 659.8
$ws.Range("J97").Value = 659.8
$ws.Range("L97").Value = 1979.4
$ws.Range("N97").Value = -2971.4
$ws.Range("H113").Value = 3133
$ws.Range("I113").Value = 3133
$ws.Range("K113").Value = 3133
$ws.Range("M113").Value = 121
$ws.Range("H132").Value = 2002.2
$ws.Range("I132").Value = 1891.3334
$ws.Range("K132").Value = 5674.0002
$ws.Range("M132").Value = -3144.0002
$ws.Range("H138").Value = 3230.5
$ws.Range("I138").Value = 1453.3889
$ws.Range("J138").Value = 4563.3335
$ws.Range("K138").Value = 4360.1667
$ws.Range("L138").Value = 13690.0005
$ws.Range("M138").Value = 779.8333000000002
$ws.Range("N138").Value = -23970.0005
$ws.Range("H140").Value = 109944
$ws.Range("J140").Value = 109944
$ws.Range("L140").Value = 109944
$ws.Range("N140").Value = -120304

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4067.6453
$ws.Range("I32").Value = 3289.5715
$ws.Range("J32").Value = 11329.667
$ws.Range("K32").Value = 3289.5715
$ws.Range("L32").Value = 11329.667
$ws.Range("M32").Value = -3002.5715
$ws.Range("N32").Value = -11903.667
$ws.Range("H74").Value = 956.8
$ws.Range("I74").Value = 886.75
$ws.Range("K74").Value = 886.75
$ws.Range("M74").Value = -12.75
$ws.Range("H77").Value = 956.8
$ws.Range("I77").Value = 886.75
$ws.Range("K77").Value = 4433.75
$ws.Range("M77").Value = -65.75
$ws.Range("H97").Value = 758.41174
$ws.Range("I97").Value = 367.85715
$ws.Range("K97").Value = 367.85715
$ws.Range("M97").Value = 128.14285
$ws.Range("H132").Value = 2233.7646
$ws.Range("I132").Value = 1853.5172
$ws.Range("J132").Value = 4439.2
$ws.Range("K132").Value = 5560.5516
$ws.Range("L132").Value = 13317.6
$ws.Range("M132").Value = -3030.5516
$ws.Range("N132").Value = -18377.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 659
$ws.Range("J22").Value = 594
$ws.Range("L22").Value = 594
$ws.Range("N22").Value = -940
$ws.Range("H134").Value = 2589.3794
$ws.Range("I134").Value = 2574.25
$ws.Range("K134").Value = 7722.75
$ws.Range("M134").Value = -5187.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8749.5
$ws.Range("I31").Value = 3000
$ws.Range("J31").Value = 10666
$ws.Range("K31").Value = 3000
$ws.Range("L31").Value = 10666
$ws.Range("M31").Value = -2705
$ws.Range("N31").Value = -11256
$ws.Range("H34").Value = 8749.5
$ws.Range("I34").Value = 3000
$ws.Range("J34").Value = 10666
$ws.Range("K34").Value = 3000
$ws.Range("L34").Value = 10666
$ws.Range("M34").Value = -2798
$ws.Range("N34").Value = -11070
$ws.Range("H58").Value = 2540
$ws.Range("I58").Value = 2399
$ws.Range("K58").Value = 2399
$ws.Range("M58").Value = -2196
$ws.Range("H74").Value = 34987.25
$ws.Range("J74").Value = 34987.25
$ws.Range("L74").Value = 34987.25
$ws.Range("N74").Value = -36735.25
$ws.Range("H77").Value = 34987.25
$ws.Range("J77").Value = 34987.25
$ws.Range("L77").Value = 104961.75
$ws.Range("N77").Value = -113697.75
$ws.Range("H122").Value = 1501.7858
$ws.Range("I122").Value = 1693.1818
$ws.Range("J122").Value = 800
$ws.Range("K122").Value = 5079.5454
$ws.Range("L122").Value = 2400
$ws.Range("M122").Value = -2629.5454
$ws.Range("N122").Value = -7300
$ws.Range("H132").Value = 2817.9062
$ws.Range("I132").Value = 2614.5386
$ws.Range("K132").Value = 7843.6158
$ws.Range("M132").Value = -5313.6158
$ws.Range("H134").Value = 4400.8335
$ws.Range("I134").Value = 4451
$ws.Range("K134").Value = 13353
$ws.Range("M134").Value = -10818
$ws.Range("H136").Value = 2540
$ws.Range("I136").Value = 2399
$ws.Range("K136").Value = 7197
$ws.Range("M136").Value = -4647

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 923.3333
$ws.Range("I5").Value = 885
$ws.Range("K5").Value = 2655
$ws.Range("M5").Value = -2543
$ws.Range("H6").Value = 99
$ws.Range("I6").Value = 99
$ws.Range("K6").Value = 297
$ws.Range("M6").Value = -184
$ws.Range("H98").Value = 4184.7144
$ws.Range("J98").Value = 3559.6
$ws.Range("L98").Value = 10678.8
$ws.Range("N98").Value = -13674.8
$ws.Range("H135").Value = 923.3333
$ws.Range("I135").Value = 885
$ws.Range("K135").Value = 7965
$ws.Range("M135").Value = -5430

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("H132").Value = 3075.0588
$ws.Range("I132").Value = 2773.9167
$ws.Range("K132").Value = 8321.750100000001
$ws.Range("M132").Value = -5791.750100000001
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2251.125
$ws.Range("I7").Value = 2251.125
$ws.Range("K7").Value = 2251.125
$ws.Range("M7").Value = -2139.125
$ws.Range("H40").Value = 4648.8335
$ws.Range("I40").Value = 4648.8335
$ws.Range("K40").Value = 4648.8335
$ws.Range("M40").Value = -4512.8335
$ws.Range("H46").Value = 349.5
$ws.Range("I46").Value = 349.5
$ws.Range("K46").Value = 349.5
$ws.Range("M46").Value = -161.5
$ws.Range("H93").Value = 2004.7
$ws.Range("I93").Value = 1956
$ws.Range("K93").Value = 1956
$ws.Range("M93").Value = -708
$ws.Range("H126").Value = 2251.125
$ws.Range("I126").Value = 2251.125
$ws.Range("K126").Value = 6753.375
$ws.Range("M126").Value = -4283.375
$ws.Range("H136").Value = 4400.8
$ws.Range("I136").Value = 4400.8
$ws.Range("K136").Value = 13202.4
$ws.Range("M136").Value = -10652.4
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 304
$ws.Range("I17").Value = 304
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 304
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -132
$ws.Range("H100").Value = 1766.3334
$ws.Range("J100").Value = 1499.5
$ws.Range("L100").Value = 2999
$ws.Range("N100").Value = -4081
$ws.Range("H132").Value = 2098.9473
$ws.Range("I132").Value = 1627.4286
$ws.Range("K132").Value = 4882.2858
$ws.Range("M132").Value = -2352.2858
$ws.Range("N17").ClearContents()

Write-Host "Applied scheduled-runner updates to 8 sheets."
